# Regenerate the "K" column (column G) values for the sale_chris 2023 sheet.
# These values represent recalculated strike-count style stats (std/mean based
# s_vals) replacing the old "Strike#" derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 7
    4  = 10
    5  = 5
    6  = 5
    7  = 7
    8  = 9
    9  = 3
    10 = 7
    11 = 6
    12 = 3
    13 = 11
    14 = 9
    15 = 10
    16 = 5
    17 = 0
    18 = 11
    19 = 8
    20 = 7
    21 = 7
    22 = 3
    23 = 2
    24 = 5
    25 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
